$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.781.81"
$ws.Range("E2").Value = "  -1.11%  "
$ws.Range("D3").Value = "1.796.36"
$ws.Range("E3").Value = "  -1.40%  "
$ws.Range("D4").Value = "'0.9998"
$ws.Range("D5").Value = "'309.29"
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "'0.4397"
$ws.Range("E7").Value = "  +4.12%  "
$ws.Range("D8").Value = "'0.3670"
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("D9").Value = "'0.07349"
$ws.Range("E9").Value = "  +1.90%  "
$ws.Range("E10").Value = "  +1.59%  "
$ws.Range("D11").Value = "'20.58"
$ws.Range("E11").Value = "  -0.98%  "
$ws.Range("D12").Value = "1.893.52"
$ws.Range("E12").Value = "  +4.35%  "
$ws.Range("D13").Value = "'6.619"
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("D14").Value = "'91.85"
$ws.Range("E14").Value = "  +2.64%  "
$ws.Range("D15").Value = "'0.07038"
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("D16").Value = "'5.258"
$ws.Range("E16").Value = "  -0.43%  "
$ws.Range("D17").Value = "'1.000"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").Value = "'0.000008642"
$ws.Range("E18").Value = "  -1.65%  "
$ws.Range("D19").Value = "'0.9996"
$ws.Range("D20").Value = "'14.72"
$ws.Range("E20").Value = "  -1.40%  "
$ws.Range("D21").Value = "26.813.32"
$ws.Range("E21").Value = "  -0.95%  "
$ws.Range("D22").Value = "'5.135"
$ws.Range("E22").Value = "  +0.29%  "
$ws.Range("E23").Value = "  -0.50%  "
$ws.Range("D24").Value = "'1.976"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "'151.45"
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").Value = "'2.197"
$ws.Range("E26").Value = "  -1.39%  "
$ws.Range("D27").Value = "'18.32"
$ws.Range("E27").Value = "  +0.39%  "
$ws.Range("D28").Value = "'5.178"
$ws.Range("E28").Value = "  -1.13%  "
$ws.Range("D29").Value = "'116.77"
$ws.Range("E29").Value = "  +0.48%  "
$ws.Range("D30").Value = "'0.08782"
$ws.Range("E30").Value = "  +0.41%  "
$ws.Range("D31").Value = "'0.7364"
$ws.Range("E31").Value = "  -0.32%  "
$ws.Range("E32").Value = "  -2.35%  "
$ws.Range("D33").Value = "'2.909"
$ws.Range("E33").Value = "  -1.20%  "
$ws.Range("D34").Value = "'4.417"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").Value = "'0.9999"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("D37").Value = "'0.01953"
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("D38").Value = "'0.05174"
$ws.Range("E38").Value = "  -1.29%  "
$ws.Range("D39").Value = "'0.5212"
$ws.Range("E39").Value = "  +3.53%  "
$ws.Range("E40").Value = "  -4.36%  "
$ws.Range("D41").Value = "'2.807"
$ws.Range("E41").Value = "  -2.16%  "
$ws.Range("D42").Value = "'0.1674"
$ws.Range("E42").Value = "  -0.65%  "
$ws.Range("D43").Value = "'0.5001"
$ws.Range("E43").Value = "  +6.09%  "
$ws.Range("D44").Value = "'8.416"
$ws.Range("E44").Value = "  -2.03%  "
$ws.Range("D45").Value = "'1.966"
$ws.Range("E45").Value = "  +4.10%  "
$ws.Range("D46").Value = "'10.28"
$ws.Range("E46").Value = "  -1.60%  "
$ws.Range("D47").Value = "'104.98"
$ws.Range("E47").Value = "  -1.19%  "
$ws.Range("D48").Value = "'0.9996"
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").Value = "'1.656"
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("D50").Value = "'0.06305"
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("D51").Value = "'0.9143"
$ws.Range("E51").Value = "  +1.49%  "

Write-Host "Updated cryptos list"